# Update countries & provincias Spain
# - Barein overtakes Kazajistan in ranking (swap display order, each keeps its own refreshed data)
# - Camerun overtakes Bosnia y Herzegovina / Afganistan / Nueva Zelanda (moves up three rows)
# - Refresh a handful of country case counts
# - Bump the "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 15:22"

# --- Alemania (row 8) ---------------------------------------------------
$ws.Range("E8").Value = 39466
$ws.Range("G8").Value = 28
$ws.Range("H8").Value = 5788

# --- Serbia (row 42) -----------------------------------------------------
$ws.Range("B42").Value = 7779
$ws.Range("C42").Value = 296
$ws.Range("D42").Value = 1152
$ws.Range("E42").Value = 6476
$ws.Range("F42").Value = 91
$ws.Range("G42").Value = 7
$ws.Range("H42").Value = 151

# --- Moldavia (row 59) ----------------------------------------------------
$ws.Range("E59").Value = 2196
$ws.Range("G59").Value = 5
$ws.Range("H59").Value = 89

# --- Barein / Kazajistan swap places (rows 62-63) --------------------------
# Barein's updated totals move it above Kazajistan.
$ws.Range("A62").Value = "Barein"
$ws.Range("B62").Value = 2588
$ws.Range("C62").Value = 70
$ws.Range("D62").Value = 1160
$ws.Range("E62").Value = 1420
$ws.Range("F62").Value = 2
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 8

$ws.Range("A63").Value = "Kazajistan"
$ws.Range("B63").Value = 2564
$ws.Range("C63").Value = 148
$ws.Range("D63").Value = 629
$ws.Range("E63").Value = 1910
$ws.Range("F63").Value = 31
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 25

# --- Islandia (row 69) -----------------------------------------------------
$ws.Range("B69").Value = 1790
$ws.Range("C69").Value = 1
$ws.Range("D69").Value = 1570
$ws.Range("E69").Value = 210
$ws.Range("F69").Value = 4

# --- Camerun moves up above Bosnia y Herzegovina / Afganistan / Nueva Zelanda
# (rows 74-77 shift down one slot each, Camerun takes row 74)
$ws.Range("A74").Value = "Camerun"
$ws.Range("B74").Value = 1518
$ws.Range("C74").Value = 88
$ws.Range("D74").Value = 697
$ws.Range("E74").Value = 768
$ws.Range("F74").Value = 28
$ws.Range("G74").Value = 10
$ws.Range("H74").Value = 53

$ws.Range("A75").Value = "Bosnia y Herzegovina"
$ws.Range("B75").Value = 1486
$ws.Range("C75").Value = 65
$ws.Range("D75").Value = 592
$ws.Range("E75").Value = 837
$ws.Range("F75").Value = 4
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 57

$ws.Range("A76").Value = "Afganistan"
$ws.Range("B76").Value = 1463
$ws.Range("C76").Value = 112
$ws.Range("D76").Value = 188
$ws.Range("E76").Value = 1228
$ws.Range("F76").Value = 7
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = 47

$ws.Range("A77").Value = "Nueva Zelanda"
$ws.Range("B77").Value = 1461
$ws.Range("C77").Value = 5
$ws.Range("D77").Value = 1118
$ws.Range("E77").Value = 325
$ws.Range("F77").Value = 1
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 18

# --- Liberia (row 143) ------------------------------------------------------
$ws.Range("B143").Value = 120
$ws.Range("C143").Value = 3
$ws.Range("G143").Value = 3
$ws.Range("H143").Value = 11
